# Auto-generated: update Pais sheet country stats per commit 'Update countries & provincias Spain'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 19:35"
$ws.Range("B4").Value = 1439715
$ws.Range("C4").Value = 9367
$ws.Range("D4").Value = 311708
$ws.Range("E4").Value = 1042239
$ws.Range("G4").Value = 571
$ws.Range("H4").Value = 85768
$ws.Range("B9").Value = 196375
$ws.Range("C9").Value = 7218
$ws.Range("E9").Value = 104396
$ws.Range("G9").Value = 397
$ws.Range("H9").Value = 13555
$ws.Range("B17").Value = 72520
$ws.Range("C17").Value = 242
$ws.Range("E17").Value = 32019
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 5337
$ws.Range("E37").Value = 6141
$ws.Range("G37").Value = 17
$ws.Range("H37").Value = 1053
$ws.Range("B42").Value = 12739
$ws.Range("C42").Value = 665
$ws.Range("D42").Value = 5676
$ws.Range("E42").Value = 6825
$ws.Range("G42").Value = 19
$ws.Range("H42").Value = 238
$ws.Range("B45").Value = 11320
$ws.Range("C45").Value = 124
$ws.Range("D45").Value = 3351
$ws.Range("E45").Value = 7547
$ws.Range("G45").Value = 13
$ws.Range("H45").Value = 422
$ws.Range("B56").Value = 6607
$ws.Range("C56").Value = 95
$ws.Range("D56").Value = 3310
$ws.Range("E56").Value = 3107
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 190
$ws.Range("D106").Value = 481
$ws.Range("E106").Value = 409
$ws.Range("B116").Value = 754
$ws.Range("C116").Value = 14
$ws.Range("D116").Value = 184
$ws.Range("E116").Value = 559
$ws.Range("A119").Value = "Georgia"
$ws.Range("B119").Value = 667
$ws.Range("C119").Value = 20
$ws.Range("D119").Value = 383
$ws.Range("E119").Value = 272
$ws.Range("F119").Value = 6
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 12
$ws.Range("A120").Value = "Zambia"
$ws.Range("B120").Value = 654
$ws.Range("C120").Value = 208
$ws.Range("D120").Value = 124
$ws.Range("E120").Value = 523
$ws.Range("F120").Value = 1
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 7
$ws.Range("A131").Value = "Republica del Chad"
$ws.Range("B131").Value = 399
$ws.Range("C131").Value = 27
$ws.Range("D131").Value = 83
$ws.Range("E131").Value = 270
$ws.Range("G131").Value = 4
$ws.Range("H131").Value = 46
$ws.Range("A132").Value = "Estado de Palestina"
$ws.Range("B132").Value = 375
$ws.Range("D132").Value = 310
$ws.Range("E132").Value = 63
$ws.Range("H132").Value = 2
$ws.Range("A144").Value = "Madagascar"
$ws.Range("B144").Value = 230
$ws.Range("C144").Value = 18
$ws.Range("D144").Value = 108
$ws.Range("E144").Value = 122
$ws.Range("F144").Value = 1
$ws.Range("H144").Value = 0
$ws.Range("A145").Value = "Santo Tome y Principe"
$ws.Range("B145").Value = 220
$ws.Range("D145").Value = 4
$ws.Range("E145").Value = 210
$ws.Range("H145").Value = 6
$ws.Range("A146").Value = "Togo"
$ws.Range("B146").Value = 219
$ws.Range("D146").Value = 96
$ws.Range("E146").Value = 112
$ws.Range("H146").Value = 11
$ws.Range("A147").Value = "Liberia"
$ws.Range("B147").Value = 215
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 105
$ws.Range("E147").Value = 90
$ws.Range("F147").Value = 0
$ws.Range("H147").Value = 20
